$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.17"
$ws.Range("E2").Value = "'0.50%"
$ws.Range("D3").Value = "'35.75"
$ws.Range("E3").Value = "'-0.51%"
$ws.Range("D4").Value = "'5.004"
$ws.Range("E4").Value = "'-1.29%"
$ws.Range("D5").Value = "'0.08067"
$ws.Range("E5").Value = "'0.04%"
$ws.Range("D6").Value = "'1.899"
$ws.Range("E6").Value = "'-2.83%"
$ws.Range("D7").Value = "'4.153"
$ws.Range("E7").Value = "'2.05%"
$ws.Range("D8").Value = "'7.885"
$ws.Range("E8").Value = "'0.87%"
$ws.Range("D9").Value = "'0.9322"
$ws.Range("E9").Value = "'0.41%"
$ws.Range("D10").Value = "'0.1241"
$ws.Range("E10").Value = "'-17.30%"
$ws.Range("D11").Value = "'0.1906"
$ws.Range("E11").Value = "'0.26%"
$ws.Range("D12").Value = "'0.09168"
$ws.Range("E12").Value = "'2.13%"
$ws.Range("D13").Value = "'0.03512"
$ws.Range("E13").Value = "'1.78%"
$ws.Range("D14").Value = "'0.09920"
$ws.Range("E14").Value = "'0.42%"
$ws.Range("D15").Value = "'0.001420"
$ws.Range("E15").Value = "'-0.24%"
$ws.Range("D16").Value = "'0.006318"
$ws.Range("E16").Value = "'7.95%"
$ws.Range("D17").Value = "'3.610"
$ws.Range("E17").Value = "'2.25%"
$ws.Range("D18").Value = "'3.115"
$ws.Range("E18").Value = "'6.38%"
$ws.Range("D19").Value = "'0.3445"
$ws.Range("E19").Value = "'-0.29%"
$ws.Range("E20").Value = "'3.08%"
$ws.Range("D21").Value = "'5.179"
$ws.Range("E21").Value = "'3.12%"
$ws.Range("E22").Value = "'5.73%"
$ws.Range("D23").Value = "'0.04421"
$ws.Range("E23").Value = "'-1.30%"
$ws.Range("D24").Value = "'0.001237"
$ws.Range("E24").Value = "'2.65%"
$ws.Range("D25").Value = "'0.004706"
$ws.Range("E25").Value = "'-2.22%"
$ws.Range("E26").Value = "'6.19%"
$ws.Range("D27").Value = "'0.0003136"
$ws.Range("E27").Value = "'4.02%"
$ws.Range("D39").Value = "'0.01944"
$ws.Range("E39").Value = "'2.46%"
$ws.Range("D40").Value = "'0.05184"
$ws.Range("E40").Value = "'8.28%"
$ws.Range("D41").Value = "'0.007558"
$ws.Range("E41").Value = "'3.17%"
$ws.Range("D42").Value = "'0.01018"
$ws.Range("E42").Value = "'-3.78%"
$ws.Range("D43").Value = "'0.1369"
$ws.Range("E43").Value = "'1.56%"
$ws.Range("D44").Value = "'0.002101"
$ws.Range("E44").Value = "'0.00%"
$ws.Range("E45").Value = "'-0.19%"
$ws.Range("D46").Value = "'0.00006366"
$ws.Range("E46").Value = "'4.22%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.39%"
$ws.Range("D48").Value = "'63.57"
$ws.Range("E48").Value = "'-1.70%"
$ws.Range("D49").Value = "'0.001663"
$ws.Range("E49").Value = "'0.20%"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.39%"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'0.39%"
